$d = $word.ActiveDocument

# 1. Grammar fix: "was downloaded" -> "were downloaded" in the Trout Lake
#    paragraph (years 2004-2013).
$rng = $d.Content.Duplicate
$rng.Find.Execute("2004-2013 was downloaded", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "2004-2013 were downloaded", 2)

# 2. Move the "_GoBack" bookmark from the "Other lakes" heading (it used to
#    sit between "O" and "ther lakes") to just after "were" in the sentence
#    we just fixed, i.e. between "were" and " downloaded".
$old = $d.Bookmarks("_GoBack")
$old.Delete()

# Re-assert the heading text in one pass so the "O" / "ther lakes" runs that
# used to straddle the bookmark collapse back into a single run.
$headingRng = $d.Content.Duplicate
$headingRng.Find.Execute("Other lakes", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "Other lakes", 2)

$target = $d.Content.Duplicate
$target.Find.Execute("2004-2013 were", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)
$bmRange = $d.Range($target.End, $target.End)
$d.Bookmarks.Add("_GoBack", $bmRange)
